$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a 5th test case (TC5) block to the "Test Suite" sheet, mirroring the
# layout used by the existing TC2/TC3/TC4 blocks (rows 15-20, 23-28, 31-36).
# The new block goes in rows 37-44: two blank rows (37, 38) followed by the
# six-row test-case block (39-44).
# ---------------------------------------------------------------------------

# Merge the description/precondition label cells first (matching the other
# blocks), *before* copying formatting over them, so the paste sets a
# uniform style across the whole merged range instead of Excel re-bordering
# the merge's inner cells.
$ws.Range("B40:D40").Merge()
$ws.Range("B41:F41").Merge()

# Copy formatting (styles/borders/fills/fonts) from the TC2 block (rows
# 15-20) onto the new block's rows (39-44) so the new rows look identical to
# the existing ones.
$ws.Range("A15:F20").Copy()
$ws.Range("A39:F44").PasteSpecial(-4122)

# Row 39: "Test Case ID:" header row
$ws.Range("A39").Value = "Test Case ID: "
$ws.Range("B39").Value = "TC5"
$ws.Range("C39").Value = "Priority (low,medium,high: "
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = "Executed by:"
$ws.Range("F39").Value = ""

# Row 40: "Description:" / "Execution Date:" row
$ws.Range("A40").Value = "Description: "
$ws.Range("B40").Value = ""
$ws.Range("C40").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = "Execution Date: "
$ws.Range("F40").Value = ""

# Row 41: "Precondition:" row
$ws.Range("A41").Value = "Precondition: "
$ws.Range("B41").Value = "O usuario devidamente autenticado e na tela inicial do sistema"
$ws.Range("C41").Value = ""
$ws.Range("D41").Value = ""
$ws.Range("E41").Value = ""
$ws.Range("F41").Value = ""

# Row 42: column headers row
$ws.Range("A42").Value = "#"
$ws.Range("B42").Value = "Steps"
$ws.Range("C42").Value = "Test Data"
$ws.Range("D42").Value = "Expected Results"
$ws.Range("E42").Value = "Execution Status (pass/fail/blocked)"
$ws.Range("F42").Value = "Actual Result"

# Row 43: step 1 (same as every other test case's first step)
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Chefe Clica para listar todas as solicitações de diárias relacionadas à sua competência."
$ws.Range("C43").Value = ""
$ws.Range("D43").Value = "SYSTEM Recupera os registros de solicitações e os exibe (em ordem decrescente pelo número da diária) em tela para o usuário."
$ws.Range("E43").Value = ""
$ws.Range("F43").Value = ""

# Row 44: step 2 (new step specific to TC5)
$ws.Range("A44").Value = 2
$ws.Range("B44").Value = "Chefe Clica para ordenar pelo nome do servidor."
$ws.Range("C44").Value = ""
$ws.Range("D44").Value = "SYSTEM Visualiza os registros de solicitações de diária ordenado pelo nome do servidor."
$ws.Range("E44").Value = ""
$ws.Range("F44").Value = ""

# Update the "Size: N test case(s))" summary cell (header block at top of
# sheet, row 3).
$ws.Range("D3").Value = "Size: 5 test case(s))"
